# The sheet contains a daily price table for "Ajo" (garlic) at the
# "Terminal La Palmera de La Serena" market. This edit adds one new
# day's record as row 108, pushing all the existing rows (108-209)
# down by one (to 109-210), which is exactly what happens when a new
# daily reading is prepended to the top of the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 108; Excel shifts rows 108:209 down to 109:210
# and the new blank row inherits the number formatting (the date style)
# from the row that used to be there.
$ws.Rows(108).Insert()

# Populate the new row 108 with the new day's data.
$ws.Cells.Item(108, 1).Value  = 8
$ws.Cells.Item(108, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(108, 3).Value  = "Coquimbo"
$ws.Cells.Item(108, 4).Value  = 44586
$ws.Cells.Item(108, 5).Value  = 4
$ws.Cells.Item(108, 6).Value  = 100112003
$ws.Cells.Item(108, 7).Value  = "Ajo"
$ws.Cells.Item(108, 8).Value  = "Chino"
$ws.Cells.Item(108, 9).Value  = "Primera"
$ws.Cells.Item(108, 10).Value = 560
$ws.Cells.Item(108, 11).Value = 19000
$ws.Cells.Item(108, 12).Value = 20000
$ws.Cells.Item(108, 13).Value = 19500
$ws.Cells.Item(108, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(108, 15).Value = "China"
$ws.Cells.Item(108, 16).Value = 1950
$ws.Cells.Item(108, 17).Value = 10
$ws.Cells.Item(108, 18).Value = "Hortaliza"
